$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# New rows 142-176: per-unit labels for displayLabel (ps/ns/us/ms), each with its own
# SingleUseId-tagged TEXT ID (col B), alignment (col D: Left then Center), LTR direction (col E).
$rows = @(
    @(142, "SingleUseId157", "displayLabel", "Left", "LTR", "ps"),
    @(143, "SingleUseId158", "displayLabel", "Left", "LTR", "ps"),
    @(144, "SingleUseId159", "displayLabel", "Left", "LTR", "ps"),
    @(145, "SingleUseId160", "displayLabel", "Left", "LTR", "ps"),
    @(146, "SingleUseId161", "displayLabel", "Left", "LTR", "ps"),
    @(147, "SingleUseId162", "displayLabel", "Left", "LTR", "ps"),
    @(148, "SingleUseId163", "displayLabel", "Left", "LTR", "ps"),
    @(149, "SingleUseId164", "displayLabel", "Center", "LTR", "ps"),
    @(150, "SingleUseId165", "displayLabel", "Center", "LTR", "ps"),
    @(151, "SingleUseId166", "displayLabel", "Center", "LTR", "ps"),
    @(152, "SingleUseId167", "displayLabel", "Center", "LTR", "ps"),
    @(153, "SingleUseId168", "displayLabel", "Center", "LTR", "ps"),
    @(154, "SingleUseId169", "displayLabel", "Center", "LTR", "ps"),
    @(155, "SingleUseId170", "displayLabel", "Center", "LTR", "ps"),
    @(156, "SingleUseId171", "displayLabel", "Center", "LTR", "ns"),
    @(157, "SingleUseId172", "displayLabel", "Center", "LTR", "ns"),
    @(158, "SingleUseId173", "displayLabel", "Center", "LTR", "ns"),
    @(159, "SingleUseId174", "displayLabel", "Center", "LTR", "ns"),
    @(160, "SingleUseId175", "displayLabel", "Center", "LTR", "ns"),
    @(161, "SingleUseId176", "displayLabel", "Center", "LTR", "ns"),
    @(162, "SingleUseId177", "displayLabel", "Center", "LTR", "ns"),
    @(163, "SingleUseId178", "displayLabel", "Center", "LTR", "us"),
    @(164, "SingleUseId179", "displayLabel", "Center", "LTR", "us"),
    @(165, "SingleUseId180", "displayLabel", "Center", "LTR", "us"),
    @(166, "SingleUseId181", "displayLabel", "Center", "LTR", "us"),
    @(167, "SingleUseId182", "displayLabel", "Center", "LTR", "us"),
    @(168, "SingleUseId183", "displayLabel", "Center", "LTR", "us"),
    @(169, "SingleUseId184", "displayLabel", "Center", "LTR", "us"),
    @(170, "SingleUseId185", "displayLabel", "Center", "LTR", "ms"),
    @(171, "SingleUseId186", "displayLabel", "Center", "LTR", "ms"),
    @(172, "SingleUseId187", "displayLabel", "Center", "LTR", "ms"),
    @(173, "SingleUseId188", "displayLabel", "Center", "LTR", "ms"),
    @(174, "SingleUseId189", "displayLabel", "Center", "LTR", "ms"),
    @(175, "SingleUseId190", "displayLabel", "Center", "LTR", "ms"),
    @(176, "SingleUseId191", "displayLabel", "Center", "LTR", "ms")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]   # B: TEXT ID
    $ws.Cells.Item($r, 3).Value = $row[2]   # C: TYPOGRAPHY NAME
    $ws.Cells.Item($r, 4).Value = $row[3]   # D: ALIGNMENT
    $ws.Cells.Item($r, 5).Value = $row[4]   # E: DIRECTION
    $ws.Cells.Item($r, 6).Value = $row[5]   # F: GB (translation text)
}
